# BangPhanCongCongViec.xlsx - update "Hoan thanh (%)" column (H) to 100%
# for all tasks, including filling in the previously-empty H22 cell, and
# move the active selection to reflect where the user ended up working.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bump completion to 100% for every task row that wasn't already there.
$ws.Range("H16").Value = 1
$ws.Range("H17").Value = 1
$ws.Range("H18").Value = 1
$ws.Range("H19").Value = 1
$ws.Range("H21").Value = 1

# Row 22 ("Thiet ke giao dien") had no completion value yet; give it the
# same percentage number formatting used by the rest of the column (H16:H21)
# before writing the value.
$h22 = $ws.Range("H22")
$h22.NumberFormat = $ws.Range("H21").NumberFormat
$h22.Font.Name = $ws.Range("H21").Font.Name
$h22.Font.Size = $ws.Range("H21").Font.Size
$h22.HorizontalAlignment = $ws.Range("H21").HorizontalAlignment
$h22.VerticalAlignment = $ws.Range("H21").VerticalAlignment
$h22.Value = 1

# Reflect the final on-screen selection/scroll position from the saved view.
$ws.Range("J24").Select()
